$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old contents (rows 4-6 need to be removed entirely)
$ws.Range("A4:A6").EntireRow.Delete()

# Update values for remaining rows
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "BUKU-BUKU TAMBAHAN T1"
$ws.Range("A3").Value = "BUKU TULIS/ALAT-ALAT"

# Set the active selection to A3, matching the target worksheet view
$ws.Range("A3").Select()
